$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 88134360
$ws.Range("B2").Value = 2071060
$ws.Range("C2").Value = "Retira - SC PALHOCA - Loja Palhoça (30)"
$ws.Range("D2").Value = "1bd"
$ws.Range("E2").Value = "Grátis"

$ws.Range("A3").Value = 88134360
$ws.Range("B3").Value = 2071060
$ws.Range("C3").Value = "Retira - SC SAO JOSE - Loja Campinas (1)"
$ws.Range("D3").Value = "1bd"
$ws.Range("E3").Value = "Grátis"

$ws.Range("A4").Value = 88134360
$ws.Range("B4").Value = 2071060
$ws.Range("C4").Value = "Entrega SC"
$ws.Range("D4").Value = "4bd"
$ws.Range("E4").Value = "R$ 29.65"
